$p = $ppt.ActivePresentation
$m = $p.SlideMaster
try {
  $tcs = $m.ThemeColorScheme
  Write-Output "tcs: $tcs count=$($tcs.Count)"
} catch {
  Write-Output "ERR tcs: $_"
}
